$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Append the new draw-result row (row 15) below the existing data, matching
# the sheet's existing convention of storing every value as literal text
# (dates, phase codes, etc. are kept as text, not "real" Excel numbers/dates).
# Columns A, C and E look like a date / number / date-time to Excel, so they
# are forced to text format ("@") before assignment to avoid automatic
# numeric/date coercion; B and D ("Pick 3", "2-1-8") are already plain text
# and don't need that treatment.

$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2025-10-01"

$ws.Range("B15").Value = "Pick 3"

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "251001"

$ws.Range("D15").Value = "2-1-8"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2025-10-01T21:38:23.744+04:00"
